# Aggiustato topAziende e inserita funzione che logga dentro la cartella dei log
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 becomes "2" (was "1"), B2 becomes "1" (was empty... well was "1" too originally
# but conceptually shifted), C2 becomes "test". Force A2/B2 to stay text (not auto-numify)
# by briefly marking the cell as Text format, writing the value, then restoring the
# default "Normal" style so no stray per-cell formatting is left behind.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "test"

# Row 3: only C3 and E3 text changes; B3/D3 stay as-is
$ws.Range("C3").Value = "test2"
$ws.Range("E3").Value = "ee"

# Row 4: B4, C4, E4 are cleared out (D4 "packaging ecologico" remains untouched)
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
